$d = $word.ActiveDocument

$replacements = @(
    @("336÷5=67, 1", "877÷4=219, 1"),
    @("821÷7=117, 2", "263÷9=29, 2"),
    @("675÷7=96, 3", "296÷9=32, 8"),
    @("149÷3=49, 2", "357÷4=89, 1"),
    @("935÷5=187, 0", "379÷8=47, 3"),
    @("885÷8=110, 5", "645÷7=92, 1"),
    @("624÷5=124, 4", "627÷4=156, 3"),
    @("658÷9=73, 1", "838÷4=209, 2"),
    @("481÷3=160, 1", "102÷4=25, 2"),
    @("742÷7=106, 0", "998÷9=110, 8"),
    @("446÷8=55, 6", "417÷6=69, 3"),
    @("843÷3=281, 0", "915÷2=457, 1"),
    @("100÷4=25, 0", "170÷9=18, 8"),
    @("438÷8=54, 6", "177÷9=19, 6"),
    @("297÷5=59, 2", "556÷6=92, 4"),
    @("641÷4=160, 1", "826÷4=206, 2"),
    @("320÷2=160, 0", "583÷7=83, 2"),
    @("610÷9=67, 7", "219÷5=43, 4"),
    @("122÷5=24, 2", "725÷6=120, 5"),
    @("153÷8=19, 1", "535÷3=178, 1"),
    @("193÷8=24, 1", "744÷4=186, 0"),
    @("343÷9=38, 1", "213÷3=71, 0"),
    @("407÷7=58, 1", "370÷5=74, 0"),
    @("381÷9=42, 3", "896÷5=179, 1"),
    @("143÷4=35, 3", "383÷4=95, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Write-Host "Replacements applied"
